$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "305.92"
Set-TextValue "E2" "0.79%"
Set-TextValue "G2" "23"
Set-TextValue "D3" "36.38"
Set-TextValue "E3" "3.21%"
Set-TextValue "G3" "23"
Set-TextValue "D4" "5.108"
Set-TextValue "E4" "0.49%"
Set-TextValue "G4" "23"
Set-TextValue "D5" "0.08082"
Set-TextValue "E5" "1.74%"
Set-TextValue "G5" "23"
Set-TextValue "D6" "1.935"
Set-TextValue "E6" "-0.79%"
Set-TextValue "G6" "23"
Set-TextValue "D7" "7.731"
Set-TextValue "E7" "-1.35%"
Set-TextValue "G7" "23"
Set-TextValue "D8" "0.9312"
Set-TextValue "E8" "0.67%"
Set-TextValue "G8" "23"
Set-TextValue "D9" "0.1459"
Set-TextValue "E9" "35.11%"
Set-TextValue "G9" "23"
Set-TextValue "D10" "0.1923"
Set-TextValue "E10" "1.77%"
Set-TextValue "G10" "23"
Set-TextValue "D11" "0.09103"
Set-TextValue "E11" "-4.21%"
Set-TextValue "G11" "23"
Set-TextValue "D12" "0.03550"
Set-TextValue "E12" "-3.15%"
Set-TextValue "G12" "23"
Set-TextValue "D13" "0.09793"
Set-TextValue "E13" "-1.30%"
Set-TextValue "G13" "23"
Set-TextValue "D14" "0.001433"
Set-TextValue "E14" "2.06%"
Set-TextValue "G14" "23"
Set-TextValue "D15" "0.005805"
Set-TextValue "E15" "-0.41%"
Set-TextValue "G15" "23"
Set-TextValue "D16" "3.523"
Set-TextValue "E16" "1.81%"
Set-TextValue "G16" "23"
Set-TextValue "D17" "4.118"
Set-TextValue "E17" "0.04%"
Set-TextValue "G17" "23"
Set-TextValue "D18" "2.932"
Set-TextValue "E18" "3.42%"
Set-TextValue "G18" "23"
Set-TextValue "D19" "0.3431"
Set-TextValue "E19" "0.42%"
Set-TextValue "G19" "23"
Set-TextValue "D20" "0.1302"
Set-TextValue "E20" "-0.82%"
Set-TextValue "G20" "23"
Set-TextValue "D21" "5.046"
Set-TextValue "G21" "23"
Set-TextValue "D22" "0.2412"
Set-TextValue "E22" "9.53%"
Set-TextValue "G22" "23"
Set-TextValue "D23" "0.04533"
Set-TextValue "E23" "0.02%"
Set-TextValue "G23" "23"
Set-TextValue "D24" "0.001213"
Set-TextValue "E24" "-1.17%"
Set-TextValue "G24" "23"
Set-TextValue "D25" "0.004851"
Set-TextValue "E25" "3.55%"
Set-TextValue "G25" "23"
Set-TextValue "D26" "0.0001248"
Set-TextValue "E26" "-0.66%"
Set-TextValue "G26" "23"
Set-TextValue "D27" "0.0004453"
Set-TextValue "E27" "-0.16%"
Set-TextValue "G27" "23"
Set-TextValue "G28" "23"
Set-TextValue "G29" "23"
Set-TextValue "G30" "23"
Set-TextValue "G31" "23"
Set-TextValue "G32" "23"
Set-TextValue "G33" "23"
Set-TextValue "G34" "23"
Set-TextValue "G35" "23"
Set-TextValue "G36" "23"
Set-TextValue "G37" "23"
Set-TextValue "G38" "23"
Set-TextValue "D39" "0.01972"
Set-TextValue "E39" "3.61%"
Set-TextValue "G39" "23"
Set-TextValue "D40" "0.04832"
Set-TextValue "E40" "2.14%"
Set-TextValue "G40" "23"
Set-TextValue "D41" "0.01119"
Set-TextValue "E41" "13.65%"
Set-TextValue "G41" "23"
Set-TextValue "D42" "0.007507"
Set-TextValue "E42" "-1.53%"
Set-TextValue "G42" "23"
Set-TextValue "D43" "0.1363"
Set-TextValue "E43" "1.27%"
Set-TextValue "G43" "23"
Set-TextValue "D44" "0.002092"
Set-TextValue "E44" "-1.54%"
Set-TextValue "G44" "23"
Set-TextValue "D45" "0.009876"
Set-TextValue "E45" "-13.24%"
Set-TextValue "G45" "23"
Set-TextValue "D46" "0.00006396"
Set-TextValue "E46" "1.75%"
Set-TextValue "G46" "23"
Set-TextValue "D47" "0.00000000752"
Set-TextValue "E47" "-0.02%"
Set-TextValue "G47" "23"
Set-TextValue "E48" "-2.70%"
Set-TextValue "G48" "23"
Set-TextValue "D49" "0.001194"
Set-TextValue "E49" "-8.44%"
Set-TextValue "G49" "23"
Set-TextValue "D50" "0.00002106"
Set-TextValue "E50" "-0.02%"
Set-TextValue "G50" "23"
Set-TextValue "D51" "0.0002005"
Set-TextValue "E51" "-0.02%"
Set-TextValue "G51" "23"
